# AMOSTO.xlsx: v1.2 switch to XMATCH
#
# 1. Lookup sheet: replace the two MATCH(...,0) lookups with XMATCH(...)
#    (XMATCH defaults to exact match, so the trailing ",0" is dropped).
# 2. ChangeLog sheet: reword the v1.1 entry describing the old MATCH fix,
#    and add a new v1.2 row describing the switch to XMATCH.
# 3. Model sheet's "Version" cell (G27) recalculates automatically from
#    MAX(ChangeLog!A4:A1011), so it will pick up the new 1.2 entry.

$wb = $excel.ActiveWorkbook

# --- Lookup sheet: MATCH() -> XMATCH() ---------------------------------
$wsLookup = $wb.Worksheets.Item("Lookup")
$wsLookup.Range("C2").Formula = "=XMATCH(B2,Parameters!A3:A5)"
$wsLookup.Range("C3").Formula = "=XMATCH(B3,Parameters!B2:C2)"
$wsLookup.Range("C4").Select()

# --- ChangeLog sheet: update notes, add v1.2 row -----------------------
$wsLog = $wb.Worksheets.Item("ChangeLog")

# Clarify the earlier fix's description (now mentions MATCH() by name and
# notes it ended with ",0").
$wsLog.Cells.Item(9, 5).Value = 'Fix lookup behavior in "MATCH()" for cattle/digestate based on user feedback (end with ,0)'

# New row documenting this release.
$wsLog.Cells.Item(10, 1).Value = 1.2
$wsLog.Cells.Item(10, 2).Value = 45033
$wsLog.Cells.Item(10, 2).NumberFormat = "d-mmm-yy"
$wsLog.Cells.Item(10, 3).Value = "AMOSTO.xlsx"
$wsLog.Cells.Item(10, 4).Value = "Sasha"
$wsLog.Cells.Item(10, 5).Value = "Switch MATCH() to XMATCH() based on MS Excel help file"

$wsLog.Range("E11").Select()

# --- Restore Model as the active sheet ---------------------------------
$wsModel = $wb.Worksheets.Item("Model")
$wsModel.Select()

"done"
